$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the relevant paragraphs by their (current) text rather than a
# hard-coded index, so the script is resilient to minor positional
# differences.
# ------------------------------------------------------------------
$firstHyperlinkIndex = 0
$researchIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($firstHyperlinkIndex -eq 0 -and $t -like "*free-texture-site.blogspot.com*") {
        $firstHyperlinkIndex = $i
    }
    if ($t -like "Research:*") {
        $researchIndex = $i
    }
}

# ------------------------------------------------------------------
# 1. Turn the empty paragraph that follows the first hyperlink into a
#    new hyperlink paragraph pointing at the skyboxes page.
# ------------------------------------------------------------------
$emptyPara = $d.Paragraphs.Item($firstHyperlinkIndex + 1)
$d.Hyperlinks.Add($emptyPara.Range, "http://www.custommapmakers.org/skyboxes.php", `
    [System.Type]::Missing, [System.Type]::Missing, `
    "http://www.custommapmakers.org/skyboxes.php") | Out-Null

# ------------------------------------------------------------------
# 2. Insert a brand new (empty) paragraph right before "Research:" and
#    relocate the existing "_GoBack" bookmark onto it (it currently
#    sits in the paragraph right after "Research:").
# ------------------------------------------------------------------
$researchPara = $d.Paragraphs.Item($researchIndex)
$researchPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($researchIndex)

# This engine (like Word itself) cannot anchor a bookmark on a fully
# empty paragraph that has no run at all, so briefly insert a
# placeholder character, wrap the bookmark around it, then delete the
# placeholder again - the bookmark collapses back down to an empty
# range inside the now-empty paragraph, exactly like the original.
$newPara.Range.InsertBefore("Z")
$placeholder = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null
$placeholder = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$placeholder.Text = ""
